# Rotate the full row contents (columns A:AY) among rows 7-12:
#   new row 7  <- old row 10
#   new row 8  <- old row 11
#   new row 9  <- old row 12
#   new row 10 <- old row 9
#   new row 11 <- old row 7
#   new row 12 <- old row 8
#
# Columns Y and AA hold date-looking text ("2018-10-03"). Excel's COM layer
# auto-coerces such strings into date serials when written back through
# Value2, which would change their stored type even though the visible
# text is identical across all six rows (so those two columns are skipped
# below - their content does not need to move, since it is already
# identical in every source/destination row).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- capture current contents of rows 7-12 (split around columns Y/AA) ---
$r7a  = $ws.Range("A7:X7").Value2
$r7b  = $ws.Range("Z7").Value2
$r7c  = $ws.Range("AB7:AY7").Value2

$r8a  = $ws.Range("A8:X8").Value2
$r8b  = $ws.Range("Z8").Value2
$r8c  = $ws.Range("AB8:AY8").Value2

$r9a  = $ws.Range("A9:X9").Value2
$r9b  = $ws.Range("Z9").Value2
$r9c  = $ws.Range("AB9:AY9").Value2

$r10a = $ws.Range("A10:X10").Value2
$r10b = $ws.Range("Z10").Value2
$r10c = $ws.Range("AB10:AY10").Value2

$r11a = $ws.Range("A11:X11").Value2
$r11b = $ws.Range("Z11").Value2
$r11c = $ws.Range("AB11:AY11").Value2

$r12a = $ws.Range("A12:X12").Value2
$r12b = $ws.Range("Z12").Value2
$r12c = $ws.Range("AB12:AY12").Value2

# --- write the rotated contents back ---
$ws.Range("A7:X7").Value2   = $r10a
$ws.Range("Z7").Value2      = $r10b
$ws.Range("AB7:AY7").Value2 = $r10c

$ws.Range("A8:X8").Value2   = $r11a
$ws.Range("Z8").Value2      = $r11b
$ws.Range("AB8:AY8").Value2 = $r11c

$ws.Range("A9:X9").Value2   = $r12a
$ws.Range("Z9").Value2      = $r12b
$ws.Range("AB9:AY9").Value2 = $r12c

$ws.Range("A10:X10").Value2   = $r9a
$ws.Range("Z10").Value2       = $r9b
$ws.Range("AB10:AY10").Value2 = $r9c

$ws.Range("A11:X11").Value2   = $r7a
$ws.Range("Z11").Value2       = $r7b
$ws.Range("AB11:AY11").Value2 = $r7c

$ws.Range("A12:X12").Value2   = $r8a
$ws.Range("Z12").Value2       = $r8b
$ws.Range("AB12:AY12").Value2 = $r8c
